$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.017.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.67%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.956.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.94%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'594.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.43%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'148.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.46%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.952.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.85%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.507"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.11%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +4.19%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +6.51%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.84%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +5.27%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'32.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.67%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.65%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.445.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.95%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.984.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.65%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.41%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.945.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.66%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'442.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.52%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.97%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D24").Value = "'11.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.80%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.93%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.78%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.40%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.93%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.37%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +16.49%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.42%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'26.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.61%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.01%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.992"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "'  +5.96%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.43%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.63%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'49.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.28%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'8.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.56%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -4.29%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.282"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.28%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'38.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -7.49%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'135.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.45%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.693.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.21%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0337"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.28%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'361.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.02%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E50").Value = "'  -0.53%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -3.44%  "
$ws.Range("E51").Style = "Normal"
